$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("REPLENISHMENT PICK")

# Add UPH header
$ws.Range("C1").Value = "UPH"

# Flip ReplenishmentPickQuantity values to positive and compute UPH
for ($r = 2; $r -le 19; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $posValue = [Math]::Abs($bCell.Value2)
    $bCell.Value = $posValue
    $ws.Cells.Item($r, 3).Value = $posValue * 15 / 38
}
